$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("D9").Value = "MBA 학생들 시험 후기 – 실험충 vs. Theorist"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/s2022-1st-term-review/#utm_source=rss&utm_medium=rss&utm_campaign=s2022-1st-term-review"

# Row 32
$ws.Range("D32").Value = "Information Gain and Mutual Information"
$ws.Range("E32").Value = "https://dodonam.tistory.com/378"

# Row 37
$ws.Range("D37").Value = "dsba_seminar"

# Row 51
$ws.Range("D51").Value = "[folium] 파이썬으로 지도 위에 마커 표시하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/folium-%ED%8C%8C%EC%9D%B4%EC%8D%AC%EC%9C%BC%EB%A1%9C-%EC%A7%80%EB%8F%84-%EC%9C%84%EC%97%90-%EB%A7%88%EC%BB%A4-%ED%91%9C%EC%8B%9C%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"

# Row 52
$ws.Range("D52").Value = "Relative Risk Regression"
